# feat: add 2022-Q1 data
#
# 1. Insert a new sheet "2022-Q1" right before the "总计" (totals) sheet,
#    using the existing "2021-Q1" sheet as a style/format template (it has
#    the same 7-column fund-holdings layout with 3 data rows that the new
#    sheet needs), then overwrite its cell values with the 2022-Q1 figures.
# 2. Prepend a new "2022-Q1" row to the "总计" sheet (shifting the existing
#    rows down by one) so the summary stays in reverse-chronological order.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Step 1: create the "2022-Q1" sheet by cloning "2021-Q1" (same layout:
# header row + 3 data rows) and placing the copy right before "总计".
# ---------------------------------------------------------------------
$template = $wb.Worksheets.Item("2021-Q1")
$total = $wb.Worksheets.Item("总计")
$template.Copy($total)

# NOTE: after Copy() shuffles sheet positions, the old $total reference can
# resolve to the freshly-inserted copy instead of "总计" -- always re-fetch
# sheet handles by name right before using them.
$new = $wb.Worksheets.Item("2021-Q1 (2)")
$new.Name = "2022-Q1"

# Clear out the template's data rows (2-4) before writing the new values,
# keeping row 1 (header) intact since the header labels only change in D1.
$new.Range("B2:H4").ClearContents()

# Header row
$new.Range("D1").Value = "基金规模"

# Row 2: 004135 申万菱信量化成长混合
$new.Range("B2").NumberFormat = "@"
$new.Range("B2").Value = "004135"
$new.Range("B2").Style = "Normal"

$new.Range("C2").NumberFormat = "@"
$new.Range("C2").Value = "申万菱信量化成长混合"
$new.Range("C2").Style = "Normal"

$new.Range("D2").NumberFormat = "@"
$new.Range("D2").Value = "0.49"
$new.Range("D2").Style = "Normal"

$new.Range("E2").NumberFormat = "@"
$new.Range("E2").Value = "86.91"
$new.Range("E2").Style = "Normal"

$new.Range("F2").NumberFormat = "@"
$new.Range("F2").Value = "2.31"
$new.Range("F2").Style = "Normal"

$new.Range("G2").NumberFormat = "@"
$new.Range("G2").Value = "0.0113"
$new.Range("G2").Style = "Normal"

$new.Range("H2").Value = 2

# Row 3: 005607 华宝中证500指数增强A
$new.Range("B3").NumberFormat = "@"
$new.Range("B3").Value = "005607"
$new.Range("B3").Style = "Normal"

$new.Range("C3").NumberFormat = "@"
$new.Range("C3").Value = "华宝中证500指数增强A"
$new.Range("C3").Style = "Normal"

$new.Range("D3").NumberFormat = "@"
$new.Range("D3").Value = "0.45"
$new.Range("D3").Style = "Normal"

$new.Range("E3").NumberFormat = "@"
$new.Range("E3").Value = "94.72"
$new.Range("E3").Style = "Normal"

$new.Range("F3").NumberFormat = "@"
$new.Range("F3").Value = "1.68"
$new.Range("F3").Style = "Normal"

$new.Range("G3").NumberFormat = "@"
$new.Range("G3").Value = "0.0076"
$new.Range("G3").Style = "Normal"

$new.Range("H3").Value = 4

# Row 4: 005608 华宝中证500指数增强C
$new.Range("B4").NumberFormat = "@"
$new.Range("B4").Value = "005608"
$new.Range("B4").Style = "Normal"

$new.Range("C4").NumberFormat = "@"
$new.Range("C4").Value = "华宝中证500指数增强C"
$new.Range("C4").Style = "Normal"

$new.Range("D4").NumberFormat = "@"
$new.Range("D4").Value = "0.23"
$new.Range("D4").Style = "Normal"

$new.Range("E4").NumberFormat = "@"
$new.Range("E4").Value = "94.72"
$new.Range("E4").Style = "Normal"

$new.Range("F4").NumberFormat = "@"
$new.Range("F4").Value = "1.68"
$new.Range("F4").Style = "Normal"

$new.Range("G4").NumberFormat = "@"
$new.Range("G4").Value = "0.0039"
$new.Range("G4").Style = "Normal"

$new.Range("H4").Value = 4

# ---------------------------------------------------------------------
# Step 2: insert a new top data row in "总计" for 2022-Q1, pushing the
# existing rows (2021-Q4, 2021-Q3, ...) down by one row.
# ---------------------------------------------------------------------
$total = $wb.Worksheets.Item("总计")
$total.Range("A2:A2").EntireRow.Insert()

# Insert drags the header row's bold/bordered format onto the new blank
# row -- strip that back to the plain "Normal" style the data rows use.
$total.Range("A2:D2").Style = "Normal"

# ... then re-apply just A2's "index column" look by copying the format
# from the cell directly below it (A3, which already carries it).
$total.Range("A3").Copy()
$total.Range("A2").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$total.Range("A2").Value = 0

$total.Range("B2").NumberFormat = "@"
$total.Range("B2").Value = "2022-Q1"
$total.Range("B2").Style = "Normal"

$total.Range("C2").Value = 3
$total.Range("D2").Value = 0.02

# Renumber the index column (A) below the inserted row: 0,1,2,3,4,5
$total.Range("A3").Value = 1
$total.Range("A4").Value = 2
$total.Range("A5").Value = 3
$total.Range("A6").Value = 4
$total.Range("A7").Value = 5
